$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Rewrite the worker account-statement table (rows 16-22):
#  - PPT/862944/JOSELYN ... rows now list periods in reverse order
#    (2401, 2312, 2311, 2310, 2309, 2308)
#  - CC/1007786694/MELISA ... row (period 2310) moves from row 19 to row 22

$rows = @(
    @{ Row = 16; Tipo = "PPT"; Doc = "862944";     Nombre = "JOSELYN ANTONIA FERNANDEZ PERNALETE"; Periodo = "2401"; Mora = 46400; Salario = 1160000 },
    @{ Row = 17; Tipo = "PPT"; Doc = "862944";     Nombre = "JOSELYN ANTONIA FERNANDEZ PERNALETE"; Periodo = "2312"; Mora = 46400; Salario = 1160000 },
    @{ Row = 18; Tipo = "PPT"; Doc = "862944";     Nombre = "JOSELYN ANTONIA FERNANDEZ PERNALETE"; Periodo = "2311"; Mora = 46400; Salario = 1160000 },
    @{ Row = 19; Tipo = "PPT"; Doc = "862944";     Nombre = "JOSELYN ANTONIA FERNANDEZ PERNALETE"; Periodo = "2310"; Mora = 46400; Salario = 1160000 },
    @{ Row = 20; Tipo = "PPT"; Doc = "862944";     Nombre = "JOSELYN ANTONIA FERNANDEZ PERNALETE"; Periodo = "2309"; Mora = 46400; Salario = 1160000 },
    @{ Row = 21; Tipo = "PPT"; Doc = "862944";     Nombre = "JOSELYN ANTONIA FERNANDEZ PERNALETE"; Periodo = "2308"; Mora = 46400; Salario = 1160000 },
    @{ Row = 22; Tipo = "CC";  Doc = "1007786694"; Nombre = "MELISA CHELIAN ALFARO PACHECO";       Periodo = "2310"; Mora = 2000;  Salario = 1500000 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 2).Value = $r.Tipo
    $ws.Cells.Item($i, 3).Value = $r.Doc
    $ws.Cells.Item($i, 4).Value = $r.Nombre
    $ws.Cells.Item($i, 5).Value = $r.Periodo
    $ws.Cells.Item($i, 6).Value = $r.Mora
    $ws.Cells.Item($i, 7).Value = $r.Salario
}
